$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.523.21'
$ws.Range("E2").Value = '  -0.62%  '
$ws.Range("D3").Value = '2.354.09'
$ws.Range("E3").Value = '  -3.00%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''541.27'
$ws.Range("E5").Value = '  +0.49%  '
$ws.Range("D6").Value = '''136.72'
$ws.Range("E6").Value = '  -4.77%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '''0.524'
$ws.Range("E8").Value = '  -9.88%  '
$ws.Range("D9").Value = '2.351.40'
$ws.Range("E9").Value = '  -3.15%  '
$ws.Range("D10").Value = '''0.104'
$ws.Range("E10").Value = '  -0.52%  '
$ws.Range("E11").Value = '  +0.30%  '
$ws.Range("D12").Value = '''5.25'
$ws.Range("E12").Value = '  -2.15%  '
$ws.Range("D13").Value = '''0.341'
$ws.Range("E13").Value = '  -1.86%  '
$ws.Range("D14").Value = '''24.44'
$ws.Range("E14").Value = '  -4.67%  '
$ws.Range("D15").Value = '2.773.99'
$ws.Range("E15").Value = '  -3.28%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '''0.0000161'
$ws.Range("E16").Value = '  -1.15%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '60.215.29'
$ws.Range("E17").Value = '  -1.01%  '
$ws.Range("D18").Value = '2.344.51'
$ws.Range("E18").Value = '  -3.48%  '
$ws.Range("D19").Value = '''10.56'
$ws.Range("E19").Value = '  -2.97%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '''4.06'
$ws.Range("E20").Value = '  -1.50%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '''312.39'
$ws.Range("E21").Value = '  -1.05%  '
$ws.Range("D22").Value = '''6.60'
$ws.Range("E22").Value = '  -4.29%  '
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("B24").Value = 'SuiNetwork'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D24").Value = '''1.89'
$ws.Range("E24").Value = '  +3.55%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '''63.07'
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("B26").Value = 'Aptos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D26").Value = '''8.58'
$ws.Range("E26").Value = '  +12.77%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").Value = '''1.00'
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").Value = '2.463.18'
$ws.Range("E28").Value = '  -3.75%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '''7.92'
$ws.Range("E29").Value = '  -2.48%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0884'
$ws.Range("E30").Value = '  -6.95%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").Value = '''1.38'
$ws.Range("E31").Value = '  -3.83%  '
$ws.Range("D32").Value = '''498.51'
$ws.Range("E32").Value = '  -3.35%  '
$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").Value = '''0.144'
$ws.Range("E33").Value = '  -1.07%  '
$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").Value = '''1.78'
$ws.Range("E34").Value = '  -4.22%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '''1.52'
$ws.Range("E35").Value = '  -2.10%  '
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").Value = '''1.00'
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").Value = '''4.57'
$ws.Range("E37").Value = '  -3.34%  '
$ws.Range("D38").Value = '''0.371'
$ws.Range("E38").Value = '  -0.42%  '
$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").Value = '''18.31'
$ws.Range("E39").Value = '  +0.41%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D40").Value = '''5.22'
$ws.Range("E40").Value = '  -6.98%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '''1.77'
$ws.Range("E41").Value = '  +2.68%  '
$ws.Range("B42").Value = 'USDe'
$ws.Range("C42").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D42").Value = '''1.00'
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").Value = '''137.01'
$ws.Range("E43").Value = '  -3.37%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = '''39.98'
$ws.Range("E44").Value = '  -0.44%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '''141.35'
$ws.Range("E45").Value = '  +1.23%  '
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value = '''2.11'
$ws.Range("E46").Value = '  -5.07%  '
$ws.Range("D47").Value = '''3.52'
$ws.Range("E47").Value = '  -1.36%  '
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").Value = '''0.0508'
$ws.Range("E48").Value = '  -3.66%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '''19.38'
$ws.Range("E49").Value = '  -7.36%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '''0.568'
$ws.Range("E50").Value = '  -1.91%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").Value = '''0.0891'
$ws.Range("E51").Value = '  -3.74%  '
